# Auto-generated Excel COM-interop script
# Applies numeric updates to the Anima_Profits workbook across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1526.0349
$ws.Range("I15").Value = 1526.0349
$ws.Range("K15").Value = 4578.1047
$ws.Range("M15").Value = -4409.1047
# row 40
$ws.Range("H40").Value = 2500
$ws.Range("J40").Value = 2600
$ws.Range("L40").Value = 2600
$ws.Range("N40").Value = -2950
# row 132
$ws.Range("H132").Value = 4207.12
$ws.Range("I132").Value = 3813.2856
$ws.Range("J132").Value = 6274.75
$ws.Range("K132").Value = 11439.8568
$ws.Range("L132").Value = 18824.25
$ws.Range("M132").Value = -8909.856800000001
$ws.Range("N132").Value = -23884.25
# row 135
$ws.Range("H135").Value = 926.94446
$ws.Range("I135").Value = 417.8
$ws.Range("J135").Value = 3472.6667
$ws.Range("K135").Value = 3760.2
$ws.Range("L135").Value = 31254.0003
$ws.Range("M135").Value = -1225.2
$ws.Range("N135").Value = -36324.0003
# row 138
$ws.Range("H138").Value = 2559.3462
$ws.Range("I138").Value = 2313.85
$ws.Range("J138").Value = 3377.6667
$ws.Range("K138").Value = 6941.549999999999
$ws.Range("L138").Value = 10133.0001
$ws.Range("M138").Value = -1801.549999999999
$ws.Range("N138").Value = -20413.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 36
$ws.Range("H36").Value = 60021.75
$ws.Range("I36").Value = 30000
$ws.Range("K36").Value = 30000
$ws.Range("M36").Value = -29654
# row 61
$ws.Range("H61").Value = 2057.2341
$ws.Range("I61").Value = 1826.8485
$ws.Range("J61").Value = 2600.2856
$ws.Range("K61").Value = 1826.8485
$ws.Range("L61").Value = 2600.2856
$ws.Range("M61").Value = -1614.8485
$ws.Range("N61").Value = -3024.2856
# row 74
$ws.Range("H74").Value = 1062.825
$ws.Range("I74").Value = 757.1070999999999
$ws.Range("J74").Value = 1776.1666
$ws.Range("K74").Value = 757.1070999999999
$ws.Range("L74").Value = 1776.1666
$ws.Range("M74").Value = 116.8929000000001
$ws.Range("N74").Value = -3524.1666
# row 77
$ws.Range("H77").Value = 1062.825
$ws.Range("I77").Value = 757.1070999999999
$ws.Range("J77").Value = 1776.1666
$ws.Range("K77").Value = 3785.5355
$ws.Range("L77").Value = 8880.833000000001
$ws.Range("M77").Value = 582.4645
$ws.Range("N77").Value = -17616.833
# row 122
$ws.Range("H122").Value = 72457.28999999999
$ws.Range("I122").Value = 111833.555
$ws.Range("J122").Value = 1580
$ws.Range("K122").Value = 335500.665
$ws.Range("L122").Value = 4740
$ws.Range("M122").Value = -333050.665
$ws.Range("N122").Value = -9640
# row 136
$ws.Range("H136").Value = 2057.2341
$ws.Range("I136").Value = 1826.8485
$ws.Range("J136").Value = 2600.2856
$ws.Range("K136").Value = 5480.5455
$ws.Range("L136").Value = 7800.8568
$ws.Range("M136").Value = -2930.5455
$ws.Range("N136").Value = -12900.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 107
$ws.Range("H107").Value = 334500
$ws.Range("I107").Value = 500250
$ws.Range("K107").Value = 500250
$ws.Range("M107").Value = -498330

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 4878.016
$ws.Range("I31").Value = 1249.1666
$ws.Range("K31").Value = 1249.1666
$ws.Range("M31").Value = -954.1666
# row 34
$ws.Range("H34").Value = 4878.016
$ws.Range("I34").Value = 1249.1666
$ws.Range("K34").Value = 1249.1666
$ws.Range("M34").Value = -1047.1666
# row 62
$ws.Range("H62").Value = 2390
$ws.Range("I62").Value = 2342.8572
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2342.8572
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1718.8572
$ws.Range("N62").Value = -3748
# row 65
$ws.Range("H65").Value = 2390
$ws.Range("I65").Value = 2342.8572
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 11714.286
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -8594.286
$ws.Range("N65").Value = -18740
# row 134
$ws.Range("H134").Value = 4991.0312
$ws.Range("I134").Value = 5661.091
$ws.Range("J134").Value = 3516.9
$ws.Range("K134").Value = 16983.273
$ws.Range("L134").Value = 10550.7
$ws.Range("M134").Value = -14448.273
$ws.Range("N134").Value = -15620.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("I5").Value = 578.1818
$ws.Range("J5").Value = 1120.625
$ws.Range("K5").Value = 1734.5454
$ws.Range("L5").Value = 3361.875
$ws.Range("M5").Value = -1622.5454
$ws.Range("N5").Value = -3585.875
# row 122
$ws.Range("H122").Value = 4545.8076
$ws.Range("J122").Value = 10208.909
$ws.Range("L122").Value = 91880.181
$ws.Range("N122").Value = -96780.181
# row 131
$ws.Range("H131").Value = 1251.3871
$ws.Range("I131").Value = 602.375
$ws.Range("J131").Value = 1477.1305
$ws.Range("K131").Value = 1807.125
$ws.Range("L131").Value = 4431.3915
$ws.Range("M131").Value = 3232.875
$ws.Range("N131").Value = -14511.3915
# row 135
$ws.Range("I135").Value = 578.1818
$ws.Range("J135").Value = 1120.625
$ws.Range("K135").Value = 5203.6362
$ws.Range("L135").Value = 10085.625
$ws.Range("M135").Value = -2668.6362
$ws.Range("N135").Value = -15155.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 2
$ws.Range("H2").Value = 224.2
$ws.Range("I2").Value = 35.8
$ws.Range("J2").Value = 412.6
$ws.Range("K2").Value = 35.8
$ws.Range("L2").Value = 412.6
$ws.Range("M2").Value = 77.2
$ws.Range("N2").Value = -638.6
# row 102
$ws.Range("H102").Value = 1610.6
$ws.Range("I102").Value = 1511.1538
$ws.Range("J102").Value = 2257
$ws.Range("K102").Value = 1511.1538
$ws.Range("L102").Value = 2257
$ws.Range("M102").Value = 110.8462
$ws.Range("N102").Value = -5501
# row 132
$ws.Range("H132").Value = 2696.9148
$ws.Range("I132").Value = 2344.4595
$ws.Range("J132").Value = 4001
$ws.Range("K132").Value = 7033.3785
$ws.Range("L132").Value = 12003
$ws.Range("M132").Value = -4503.3785
$ws.Range("N132").Value = -17063

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 4520.837
$ws.Range("I22").Value = 2260.606
$ws.Range("J22").Value = 9182.5625
$ws.Range("K22").Value = 2260.606
$ws.Range("L22").Value = 9182.5625
$ws.Range("M22").Value = -1965.606
$ws.Range("N22").Value = -9772.5625
# row 27
$ws.Range("H27").Value = 4520.837
$ws.Range("I27").Value = 2260.606
$ws.Range("J27").Value = 9182.5625
$ws.Range("K27").Value = 2260.606
$ws.Range("L27").Value = 9182.5625
$ws.Range("M27").Value = -2153.606
$ws.Range("N27").Value = -9396.5625
# row 46
$ws.Range("H46").Value = 2724.25
$ws.Range("I46").Value = 2724.25
$ws.Range("K46").Value = 2724.25
$ws.Range("M46").Value = -2536.25
# row 68
$ws.Range("H68").Value = 1569.2307
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 1600
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 1600
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3098
# row 71
$ws.Range("H71").Value = 1569.2307
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 1600
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 8000
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -15488
# row 132
$ws.Range("H132").Value = 1607.9807
$ws.Range("I132").Value = 1156.8572
$ws.Range("J132").Value = 3502.7
$ws.Range("K132").Value = 3470.5716
$ws.Range("L132").Value = 10508.1
$ws.Range("M132").Value = -940.5715999999998
$ws.Range("N132").Value = -15568.1
# row 136
$ws.Range("H136").Value = 4903164.5
$ws.Range("I136").Value = 1030.1428
$ws.Range("K136").Value = 3090.4284
$ws.Range("M136").Value = -540.4284000000002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# row 126
$ws.Range("H126").Value = 3700.8572
$ws.Range("I126").Value = 2064
$ws.Range("J126").Value = 5883.3335
$ws.Range("K126").Value = 6192
$ws.Range("L126").Value = 17650.0005
$ws.Range("M126").Value = -3722
$ws.Range("N126").Value = -22590.0005
# row 132
$ws.Range("H132").Value = 5378941
$ws.Range("I132").Value = 2821.5264
$ws.Range("J132").Value = 13891130
$ws.Range("K132").Value = 8464.5792
$ws.Range("L132").Value = 41673390
$ws.Range("M132").Value = -5934.5792
$ws.Range("N132").Value = -41678450

